# outputLayout.xlsx update:
#  - Insert a new worksheet "02_Law of Indices" between "01_Change of
#    Subjects" and "20_Properties of Circles", populate it with its
#    question-code table, and make it the active sheet/selection.
#  - Update the selection on "01_Change of Subjects" to span its full
#    used range (A1:B13) instead of the old single-cell selection, and
#    hand off the active-tab/tabSelected flag to the new sheet.

$wb = $excel.ActiveWorkbook

# --- 01_Change of Subjects: refresh its selection before we move focus ---
$wsSubjects = $wb.Worksheets.Item(1)
[void]$wsSubjects.Range("A1:B13").Select()

# --- insert "02_Law of Indices" right after "01_Change of Subjects" ---
$wsIndices = $wb.Worksheets.Add($null, $wsSubjects)
$wsIndices.Name = "02_Law of Indices"

$rows = @(
  @(2,  "Question Code"),
  @(1,  "DSE12PII_Q01"),
  @(2,  "DSE14PII_Q01"),
  @(3,  "DSE15PII_Q02"),
  @(4,  "DSE19PII_Q02"),
  @(5,  "DSE22PII_Q02"),
  @(6,  "DSE24PII_Q02"),
  @(7,  "DSE20PII_Q01"),
  @(8,  "DSE21PII_Q01"),
  @(9,  "DSE23PII_Q03"),
  @(10, "DSE17PII_Q02"),
  @(11, "DSE18PII_Q01"),
  @(12, "DSE13PII_Q01"),
  @(13, "DSE16PII_Q01"),
  @(14, "DSESPPII_Q01")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $wsIndices.Cells.Item($r, 1).Value = $rows[$i][0]
    $wsIndices.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Matches the committed file: new sheet becomes the active tab, with
# A11:A15 selected (anchored at A11).
[void]$wsIndices.Range("A11:A15").Select()
